# #5: property boat&car done
#
# The car ("汽車") sheet -- the 2nd worksheet in the workbook -- gets a new
# "capacity" column plus the usual trailing metadata columns
# (property_category, category, date, legislator_name, legislator_id,
# source_file, index) that the other property sheets already have.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# ---- Header row (row 1) ------------------------------------------------
$ws.Range("B1").Value = "name"
$ws.Range("C1").Value = "capacity"
$ws.Range("D1").Value = "owner"
$ws.Range("E1").Value = "register_date"
$ws.Range("F1").Value = "register_reason"
$ws.Range("G1").Value = "acquire_value"
$ws.Range("H1").Value = "property_category"
$ws.Range("I1").Value = "category"
$ws.Range("J1").Value = "date"
$ws.Range("K1").Value = "legislator_name"
$ws.Range("L1").Value = "legislator_id"
$ws.Range("M1").Value = "source_file"
$ws.Range("N1").Value = "index"

# Match the bold/centered/bordered "header" look already used for B1:G1
$headerRange = $ws.Range("B1:N1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

# ---- Data row (row 2) ---------------------------------------------------
$ws.Range("A2").Value = 40
$ws.Range("B2").Value = "LEXUSRX350"
$ws.Range("C2").Value = 3456
$ws.Range("D2").Value = "張嘉郡"
$ws.Range("E2").Value = "100年04月08曰"
$ws.Range("F2").Value = "買賣"
$ws.Range("G2").Value = 2400000
$ws.Range("H2").Value = "land"
$ws.Range("I2").Value = "normal"

# Keep the date looking like plain text ("2012-04-20") instead of letting
# it get auto-converted into a date serial number.
$ws.Range("J2").NumberFormat = "@"
$ws.Range("J2").Value = "2012-04-20"

$ws.Range("K2").Value = "張嘉郡"
$ws.Range("L2").Value = 1719
$ws.Range("M2").Value = "tmp1fff1"
$ws.Range("N2").Value = 40
